# 📊 Horarios actualizados Línea 141 - 1106
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with the
# latest scrape (Última actualización: 03:52:29).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 03:52:29"
$ws1.Range("A3").Value = "Total filas: 7"

# Existing rows 6-10 shift to the new scraped order/content.
$ws1.Range("A6").Value = "03:52:29"
$ws1.Range("B6").Value = "04:01"
$ws1.Range("C6").Value = "81_EL PELIGRO"
$ws1.Range("D6").Value = 9
$ws1.Range("E6").Value = "LP1912"

$ws1.Range("A7").Value = "03:52:29"
$ws1.Range("B7").Value = "04:46"
$ws1.Range("C7").Value = "215A_EL PATO"
$ws1.Range("D7").Value = 54
$ws1.Range("E7").Value = "LP1912"

$ws1.Range("A8").Value = "03:52:29"
$ws1.Range("B8").Value = "04:53"
$ws1.Range("C8").Value = "11_ETCHEVERRY"
$ws1.Range("D8").Value = 61
$ws1.Range("E8").Value = "LP1912"

$ws1.Range("A9").Value = "03:52:29"
$ws1.Range("B9").Value = "05:16"
$ws1.Range("C9").Value = "17_ROMERO"
$ws1.Range("D9").Value = 84
$ws1.Range("E9").Value = "LP1912"

$ws1.Range("A10").Value = "03:52:29"
$ws1.Range("B10").Value = "05:22"
$ws1.Range("C10").Value = "23_HERNANDEZ"
$ws1.Range("D10").Value = 90
$ws1.Range("E10").Value = "LP1912"

# New rows 11-12.
$ws1.Range("A11").Value = "03:52:29"
$ws1.Range("B11").Value = "05:34"
$ws1.Range("C11").Value = "215B_EL PATO"
$ws1.Range("D11").Value = 102
$ws1.Range("E11").Value = "LP1912"

$ws1.Range("A12").Value = "03:52:29"
$ws1.Range("B12").Value = "05:46"
$ws1.Range("C12").Value = "15_ABASTO"
$ws1.Range("D12").Value = 114
$ws1.Range("E12").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 03:52:29"

$ws2.Range("A6").Value = "03:52:29"
$ws2.Range("B6").Value = "04:46"
$ws2.Range("C6").Value = "215A_EL PATO"
$ws2.Range("D6").Value = 54
$ws2.Range("E6").Value = "LP1912"

$ws2.Range("A7").Value = "03:52:29"
$ws2.Range("B7").Value = "05:34"
$ws2.Range("C7").Value = "215B_EL PATO"
$ws2.Range("D7").Value = 102
$ws2.Range("E7").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 03:52:29"
$ws3.Range("A3").Value = "Total filas: 1"

$ws3.Range("A5").Value = "Hora_Scrap"
$ws3.Range("B5").Value = "Hora_Llegada"
$ws3.Range("C5").Value = "Linea"
$ws3.Range("D5").Value = "Minutos"
$ws3.Range("E5").Value = "Parada"

$ws3.Range("A6").Value = "03:52:29"
$ws3.Range("B6").Value = "05:44"
$ws3.Range("C6").Value = "215A_LA PLATA"
$ws3.Range("D6").Value = 112
$ws3.Range("E6").Value = "L6173"
